$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "#09-파이썬(Python) 함수, lambda 함수, 인수(arguments), 매개변수(parameters)"
$ws.Range("E4").Value = "https://teddylee777.github.io/python/python-tutorial-09"

$ws.Range("D5").Value = "고유함수 전개"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/06/01/eigenfunction_expansions.html"

$ws.Range("D9").Value = "[공지] 6월 1일 2차 설명회"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/notice-webinar-20210601/#utm_source=rss&utm_medium=rss&utm_campaign=notice-webinar-20210601"

$ws.Range("D16").Value = "SS-CAM: Smoothed Score-CAM for Sharper Visual Feature Localization 내용 정리 [XAI-10]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/152"

$ws.Range("D21").Value = "[제테크] 6월 공모주(1) 엘비루셈"
$ws.Range("E21").Value = "https://ms-review.tistory.com/16"

$ws.Range("D42").Value = "GetPrivateProfileString, GetPrivateProfileInt 오류"
$ws.Range("E42").Value = "https://kjk92.tistory.com/70"

$ws.Range("D51").Value = "[git] 기존 원격저장소와의 연결 끊고 새로운 원격저장소와 연결하기"
$ws.Range("E51").Value = "https://bskyvision.com/1206"
